# Update investment capacity results for sheets "2025", "2030", "2035"
# with fresh values received from server (row 2 of each sheet).

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "A2" = 0.1851651463765911
        "B2" = 0.02589925663510731
        "E2" = 0.1734801324400598
        "G2" = 0.1304011109248329
        "H2" = 0.4039567342501246
        "I2" = 0.9957251134364422
        "N2" = 8.345858816770113
        "O2" = 5.575886644528953
    }
    "2030" = @{
        "A2" = 0.1248156211833492
        "B2" = 0.140790332440046
        "E2" = 0.08898626755994007
        "G2" = 0
        "H2" = 0
        "I2" = 0.5294595865635581
        "M2" = 0
        "N2" = 4.632083005978622
        "O2" = 2.917558375056132
    }
    "2035" = @{
        "A2" = 0.3547829215152175
        "B2" = 0.02485964336490662
        "E2" = 0
        "I2" = 0.490042408230958
        "M2" = 0.02580438963328352
        "N2" = 4.044634463648775
        "O2" = 5.839494646039698
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $cellRefStr = [string]$cellRef
        $ws.Range($cellRefStr).Value = $cellValues[$cellRef]
    }
}
